# Applies the diff: updates the date line and the 25 division-problem
# table cells. Some new values coincide with other cells' old values
# (e.g. '83÷5=16, 3' is simultaneously an old value and a new value),
# so replacements are ordered to avoid a later find-and-replace
# accidentally matching text that an earlier step just produced.
$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-24 Thursday", "2024-10-25 Friday"),
    @("61÷9=6, 7", "49÷3=16, 1"),
    @("39÷8=4, 7", "43÷4=10, 3"),
    @("47÷7=6, 5", "92÷2=46, 0"),
    @("10÷7=1, 3", "70÷8=8, 6"),
    @("85÷3=28, 1", "76÷4=19, 0"),
    @("93÷2=46, 1", "20÷9=2, 2"),
    @("92÷9=10, 2", "94÷4=23, 2"),
    @("20÷8=2, 4", "71÷6=11, 5"),
    @("10÷9=1, 1", "13÷4=3, 1"),
    @("51÷2=25, 1", "53÷6=8, 5"),
    @("47÷3=15, 2", "79÷7=11, 2"),
    @("60÷9=6, 6", "19÷6=3, 1"),
    @("57÷4=14, 1", "49÷6=8, 1"),
    @("43÷5=8, 3", "94÷9=10, 4"),
    @("12÷2=6, 0", "81÷5=16, 1"),
    @("18÷3=6, 0", "75÷6=12, 3"),
    @("77÷7=11, 0", "15÷7=2, 1"),
    @("36÷7=5, 1", "87÷5=17, 2"),
    @("83÷5=16, 3", "18÷4=4, 2"),
    @("49÷5=9, 4", "83÷5=16, 3"),
    @("89÷3=29, 2", "42÷8=5, 2"),
    @("50÷2=25, 0", "38÷3=12, 2"),
    @("63÷8=7, 7", "55÷5=11, 0"),
    @("20÷5=4, 0", "70÷8=8, 6"),
    @("73÷4=18, 1", "73÷5=14, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found -> $old"
    }
}

Write-Output "Replacements applied: $($replacements.Count)"
